# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the 15a85270-... file (row 2) on the zh-cn and de-de status sheets,
# and roll the corresponding "Latest HO Xliff Generate Date" forward on
# the Overview sheet.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: new handoff / handback timestamps for the 15a85270 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 10:56:41"
$wsZhCn.Range("K2").Value = "2016-08-28 10:57:10"

# de-de sheet: new handoff / handback timestamps for the 15a85270 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 10:56:46"
$wsDeDe.Range("K2").Value = "2016-08-28 10:57:18"

# Overview sheet: Latest HO Xliff Generate Date for the 15a85270 row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 10:56:46"
